$wb = $excel.ActiveWorkbook

# --- Rename existing "project_data" sheet to "project_data_single" ---
$wsSingle = $wb.Worksheets.Item("project_data")
$wsSingle.Name = "project_data_single"

# --- Add the new "project_data_multiple" sheet right after it ---
$wsMultiple = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSingle)
$wsMultiple.Name = "project_data_multiple"

# Populate it with the same header/first data row as project_data_single,
# plus one additional project row.
$wsMultiple.Range("A1").Value = "project_id"
$wsMultiple.Range("B1").Value = "project_name"
$wsMultiple.Range("A1:B1").Font.Bold = $true

$wsMultiple.Range("A2").Value = 141
$wsMultiple.Range("B2").Value = "Certainly a Project GmbH & Co. KG"

$wsMultiple.Range("A3").Value = 178
$wsMultiple.Range("B3").Value = "Another Project GmbH"

# Selection on the new sheet ends up on the row just after the data.
$wsMultiple.Range("A4").Select()

# --- Update the selection on project_data_single ---
$wsSingle.Range("A1:B2").Select()

# project_data_multiple is the sheet left active/in front when the file is saved.
$wsMultiple.Activate()
